# Update the "last updated" timestamp string
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 18 de Agosto de 2020 a las 09:13"

# Swap the two rows for "Islas Malvinas" (row 213) and "Montserrat" (row 214)
# so that Montserrat now sorts above Islas Malvinas.
$ws.Range("A213").Value = "Montserrat"
$ws.Range("B213").Value = 13
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 12
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 13
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0

# Update COVID statistics for several countries (row numbers identified by
# their "Pais" (A) column matching the country name).

# Kirguistan (row 56)
$ws.Range("B56").Value = 42146
$ws.Range("C56").Value = 155
$ws.Range("D56").Value = 34855
$ws.Range("E56").Value = 5793
$ws.Range("G56").Value = 2
$ws.Range("H56").Value = 1498

# Armenia (row 57)
$ws.Range("B57").Value = 41846
$ws.Range("C57").Value = 145
$ws.Range("D57").Value = 34982
$ws.Range("E57").Value = 6032
$ws.Range("G57").Value = 8
$ws.Range("H57").Value = 832

# Hungria (row 108)
$ws.Range("B108").Value = 4970
$ws.Range("C108").Value = 24
$ws.Range("D108").Value = 3631
$ws.Range("E108").Value = 730
$ws.Range("G108").Value = 1
$ws.Range("H108").Value = 609

# Lituania (row 129)
$ws.Range("B129").Value = 2474
$ws.Range("C129").Value = 38
$ws.Range("D129").Value = 1733
$ws.Range("E129").Value = 660

# Georgia (row 147)
$ws.Range("B147").Value = 1351
$ws.Range("C147").Value = 10
$ws.Range("E147").Value = 242

# Letonia (row 149)
$ws.Range("D149").Value = 1093
$ws.Range("E149").Value = 197
$ws.Range("G149").Value = 1
$ws.Range("H149").Value = 33
